# Generate Report for Handback
# Updates the localization-status report after a handback: the status
# moves from "Ready for handoff" to "Handed back: in sync with en-US",
# the handback timestamps are refreshed, and any stale handback-version
# error details are cleared now that the target files are back in sync.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# -- Status column (shared across Overview + each language sheet) --------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value     = $newStatus
$dede.Range("C2").Value     = $newStatus

# -- Refresh "Latest Handback DateTime" for each language -----------------
$zhcn.Range("K2").Value = "2016-09-03 14:53:31"
$dede.Range("K2").Value = "2016-09-03 14:53:38"

# -- Clear the stale "Error Detail" now that handback is in sync ----------
$zhcn.Range("P2").Value = " "
$zhcn.Range("P2").ClearContents()
$dede.Range("P2").Value = " "
$dede.Range("P2").ClearContents()

# -- Widen the Status columns / shrink the now-empty Error Detail column --
$overview.Columns.Item(5).ColumnWidth  = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth  = 29.9777047293527
$zhcn.Columns.Item(3).ColumnWidth      = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth     = 13.7470528738839
$dede.Columns.Item(3).ColumnWidth      = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth     = 13.7470528738839
